# Fill previously-empty Column C (GARANTI) cells on the BENCHMARK sheet
# with benchmark values, matching column D/E where applicable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

$ws.Range("C3").Value  = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("C4").Value  = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("C5").Value  = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("C6").Value  = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("C8").Value  = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("C9").Value  = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("C10").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("C11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("C12").Value = "WU: 1.000,01 USD–9,51 USD"
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("C14").Value = "40.000 TL - 1.904,76 TL"
